$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 1: string value updates (reordered/extended shared strings) ---
$ws.Range("A1").Value = "negative"
$ws.Range("J1").Value = "positive"

# Row 2 (column headers) is unchanged by this edit - leave as-is.

# --- Add new row 30 by copying formatting (styles) from row 29, then fill values ---
$ws.Range("A29:H29").Copy($ws.Range("A30:H30"))

# --- Row-by-row data updates (rows 3-30) ---
# row 3
$ws.Range("A3").Value = "poorly"
$ws.Range("B3").Value = 0.9565217391304348
$ws.Range("C3").Value = 44
$ws.Range("D3").Value = 44
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 2
$ws.Range("J3").Value = "wonderful"
$ws.Range("K3").Value = 0.875
$ws.Range("L3").Value = 49
$ws.Range("M3").Value = 49
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 7

# row 4
$ws.Range("A4").Value = "disappointing"
$ws.Range("B4").Value = 0.8409090909090909
$ws.Range("C4").Value = 37
$ws.Range("D4").Value = 37
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 7
$ws.Range("J4").Value = "awesome"
$ws.Range("K4").Value = 0.8461538461538461
$ws.Range("L4").Value = 55
$ws.Range("M4").Value = 55
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 10

# row 5
$ws.Range("A5").Value = "returned"
$ws.Range("B5").Value = 0.7894736842105263
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 30
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 8
$ws.Range("J5").Value = "favorite"
$ws.Range("K5").Value = 0.7311827956989247
$ws.Range("L5").Value = 68
$ws.Range("M5").Value = 68
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 25

# row 6
$ws.Range("A6").Value = "however"
$ws.Range("B6").Value = 0.78125
$ws.Range("C6").Value = 50
$ws.Range("D6").Value = 50
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 14
$ws.Range("J6").Value = "classic"
$ws.Range("K6").Value = 0.6037735849056604
$ws.Range("L6").Value = 32
$ws.Range("M6").Value = 32
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 21

# row 7
$ws.Range("A7").Value = "broke"
$ws.Range("B7").Value = 0.7572815533980582
$ws.Range("C7").Value = 156
$ws.Range("D7").Value = 156
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = 50
$ws.Range("J7").Value = "excellent"
$ws.Range("K7").Value = 0.5
$ws.Range("L7").Value = 32
$ws.Range("M7").Value = 32
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 32

# row 8
$ws.Range("A8").Value = "disappointed"
$ws.Range("B8").Value = 0.7204301075268817
$ws.Range("C8").Value = 134
$ws.Range("D8").Value = 134
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = 52
$ws.Range("J8").Value = "great"
$ws.Range("K8").Value = 0.3565573770491803
$ws.Range("L8").Value = 435
$ws.Range("M8").Value = 435
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 785

# row 9
$ws.Range("A9").Value = "poor"
$ws.Range("B9").Value = 0.7183098591549296
$ws.Range("C9").Value = 51
$ws.Range("D9").Value = 51
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = $false
$ws.Range("H9").Value = 20
$ws.Range("J9").Value = "love"
$ws.Range("K9").Value = 0.3113342898134864
$ws.Range("L9").Value = 217
$ws.Range("M9").Value = 217
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 480

# row 10
$ws.Range("A10").Value = "waste"
$ws.Range("B10").Value = 0.6486486486486487
$ws.Range("C10").Value = 96
$ws.Range("D10").Value = 96
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = $false
$ws.Range("H10").Value = 52
$ws.Range("J10").Value = "loves"
$ws.Range("K10").Value = 0.2448132780082987
$ws.Range("L10").Value = 118
$ws.Range("M10").Value = 118
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 364

# row 11
$ws.Range("A11").Value = "smaller"
$ws.Range("B11").Value = 0.5798319327731093
$ws.Range("C11").Value = 69
$ws.Range("D11").Value = 69
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = $false
$ws.Range("H11").Value = 50
$ws.Range("J11").Value = "best"
$ws.Range("K11").Value = 0.2416666666666667
$ws.Range("L11").Value = 29
$ws.Range("M11").Value = 29
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 91

# row 12
$ws.Range("A12").Value = "junk"
$ws.Range("B12").Value = 0.5636363636363636
$ws.Range("C12").Value = 31
$ws.Range("D12").Value = 31
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = $false
$ws.Range("H12").Value = 24
$ws.Range("J12").Value = "perfect"
$ws.Range("K12").Value = 0.2228915662650602
$ws.Range("L12").Value = 37
$ws.Range("M12").Value = 37
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 129

# row 13
$ws.Range("A13").Value = "small"
$ws.Range("B13").Value = 0.4782608695652174
$ws.Range("C13").Value = 165
$ws.Range("D13").Value = 165
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = $false
$ws.Range("H13").Value = 180
$ws.Range("J13").Value = "loved"
$ws.Range("K13").Value = 0.1957186544342508
$ws.Range("L13").Value = 64
$ws.Range("M13").Value = 64
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 263

# row 14
$ws.Range("A14").Value = "broken"
$ws.Range("B14").Value = 0.4698795180722892
$ws.Range("C14").Value = 39
$ws.Range("D14").Value = 39
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = $false
$ws.Range("H14").Value = 44
$ws.Range("J14").Value = "fun"
$ws.Range("K14").Value = 0.07887817703768624
$ws.Range("L14").Value = 90
$ws.Range("M14").Value = 90
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 1051

# row 15
$ws.Range("A15").Value = "plastic"
$ws.Range("B15").Value = 0.4330708661417323
$ws.Range("C15").Value = 55
$ws.Range("D15").Value = 55
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = $false
$ws.Range("H15").Value = 72
$ws.Range("J15").Value = "game"
$ws.Range("K15").Value = 0.03376623376623376
$ws.Range("L15").Value = 52
$ws.Range("M15").Value = 53
$ws.Range("N15").Value = 0.98
$ws.Range("O15").Value = 0.02000000000000002
$ws.Range("P15").Value = $true
$ws.Range("Q15").Value = 1488

# row 16
$ws.Range("A16").Value = "apart"
$ws.Range("B16").Value = 0.4210526315789473
$ws.Range("C16").Value = 40
$ws.Range("D16").Value = 40
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = $false
$ws.Range("H16").Value = 55

# row 17
$ws.Range("A17").Value = "difficult"
$ws.Range("B17").Value = 0.3370786516853932
$ws.Range("C17").Value = 30
$ws.Range("D17").Value = 30
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = $false
$ws.Range("H17").Value = 59

# row 18
$ws.Range("A18").Value = "ok"
$ws.Range("B18").Value = 0.3125
$ws.Range("C18").Value = 40
$ws.Range("D18").Value = 40
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = $false
$ws.Range("H18").Value = 88

# row 19
$ws.Range("A19").Value = "thought"
$ws.Range("B19").Value = 0.3118811881188119
$ws.Range("C19").Value = 63
$ws.Range("D19").Value = 63
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = $false
$ws.Range("H19").Value = 139

# row 20
$ws.Range("A20").Value = "cheap"
$ws.Range("B20").Value = 0.2843601895734597
$ws.Range("C20").Value = 60
$ws.Range("D20").Value = 60
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = $false
$ws.Range("H20").Value = 151

# row 21
$ws.Range("A21").Value = "size"
$ws.Range("B21").Value = 0.2061855670103093
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = $false
$ws.Range("H21").Value = 154

# row 22
$ws.Range("A22").Value = "item"
$ws.Range("B22").Value = 0.1811594202898551
$ws.Range("C22").Value = 50
$ws.Range("D22").Value = 50
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = $false
$ws.Range("H22").Value = 226

# row 23
$ws.Range("A23").Value = "money"
$ws.Range("B23").Value = 0.180379746835443
$ws.Range("C23").Value = 57
$ws.Range("D23").Value = 57
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = $false
$ws.Range("H23").Value = 259

# row 24
$ws.Range("A24").Value = "work"
$ws.Range("B24").Value = 0.1772151898734177
$ws.Range("C24").Value = 56
$ws.Range("D24").Value = 56
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = $false
$ws.Range("H24").Value = 260

# row 25
$ws.Range("A25").Value = "would"
$ws.Range("B25").Value = 0.1750741839762611
$ws.Range("C25").Value = 118
$ws.Range("D25").Value = 118
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = $false
$ws.Range("H25").Value = 556

# row 26
$ws.Range("A26").Value = "better"
$ws.Range("B26").Value = 0.1355140186915888
$ws.Range("C26").Value = 29
$ws.Range("D26").Value = 29
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = $false
$ws.Range("H26").Value = 185

# row 27
$ws.Range("A27").Value = "product"
$ws.Range("B27").Value = 0.1277533039647577
$ws.Range("C27").Value = 58
$ws.Range("D27").Value = 58
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = $false
$ws.Range("H27").Value = 396

# row 28
$ws.Range("A28").Value = "price"
$ws.Range("B28").Value = 0.1235632183908046
$ws.Range("C28").Value = 43
$ws.Range("D28").Value = 43
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = $false
$ws.Range("H28").Value = 305

# row 29
$ws.Range("A29").Value = "use"
$ws.Range("B29").Value = 0.1013698630136986
$ws.Range("C29").Value = 37
$ws.Range("D29").Value = 37
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = $false
$ws.Range("H29").Value = 328

# row 30
$ws.Range("A30").Value = "like"
$ws.Range("B30").Value = 0.0625
$ws.Range("C30").Value = 38
$ws.Range("D30").Value = 38
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = $false
$ws.Range("H30").Value = 570
